# Update latest output (run 42)
$wb = $excel.ActiveWorkbook

# ---- Sheet "Schedule": recompute cost / unit-cost for rows 3 & 4 ----
$sched = $wb.Worksheets.Item("Schedule")
$sched.Range("E3").Value = 352.1131185
$sched.Range("F3").Value = 23.28790466269841
$sched.Range("E4").Value = -149.15047875
$sched.Range("F4").Value = -3.288149884259259

# ---- Sheet "Detailed": refreshed Price column (and a couple of Type flips) ----
$det = $wb.Worksheets.Item("Detailed")

$det.Range("B37").Value = -6.80121
$det.Range("B38").Value = -6

$det.Range("B39").Value = -3.1159
$det.Range("C39").Value = "historical"

$det.Range("B40").Value = 34.26695
$det.Range("C40").Value = "historical"

$det.Range("B41").Value = 36.25056
$det.Range("B42").Value = 35.87161
$det.Range("B43").Value = 20.51366
$det.Range("B44").Value = 17.88508
$det.Range("B45").Value = 62.33685
$det.Range("B46").Value = 62.33685
$det.Range("B47").Value = 61.94424
$det.Range("B48").Value = 56.98
$det.Range("B49").Value = 43.27337

$det.Range("B52").Value = 41.38585
$det.Range("B53").Value = 40.54

$det.Range("B59").Value = 59.08177

$det.Range("B62").Value = 56.98

$det.Range("B66").Value = -0.88256
$det.Range("B67").Value = -5.11737
$det.Range("B68").Value = 0
$det.Range("B69").Value = -5.2121
$det.Range("B70").Value = -5.51
$det.Range("B71").Value = -5.51
$det.Range("B72").Value = -5.64248
$det.Range("B73").Value = -5.01
$det.Range("B74").Value = -1.40538
$det.Range("B75").Value = -2.67373

$det.Range("B77").Value = -9.99
$det.Range("B78").Value = -12.3505

$det.Range("B80").Value = -17.41389
$det.Range("B81").Value = -13.9999
$det.Range("B82").Value = -6.41446
$det.Range("B83").Value = -11
$det.Range("B84").Value = -9.5
$det.Range("B85").Value = -5.17024
$det.Range("B86").Value = -6.19141
$det.Range("B87").Value = -5.9299

$det.Range("B90").Value = 13.59537
$det.Range("B91").Value = 9.56921
$det.Range("B92").Value = 36.0601
